# Update the build timestamp embedded in version strings throughout the
# workbook: "February 03 2026 17.29.55 EST" -> "February 03 2026 18.05.36 EST"

$wb = $excel.ActiveWorkbook

$newStamp = "February 03 2026 18.05.36 EST"

# --- "About" sheet -------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A2").Value = "Version: Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newStamp)"

$wsAbout.Range("A6").Value = "Recommended Citation:  ""Global Energy Monitor, Coal mine boundaries and methane sources for Sima Coal Mine, China, M1203, version 'Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newStamp)'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet -------------------------------
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 8; $row++) {
    $wsData.Cells.Item($row, 19).Value = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newStamp)"
}
